$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1922.932062252702
$ws.Range("D2").Value = 1191.040285844521
